# Game of Life instructions: resize the state-chart drawing and add a new
# dashed "Straight Connector 2" inside the locked-canvas drawing, renumbering
# the drawing's docPr id from 2 to 1 (mirrors the target commit's XML diff).

$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the (only) inline drawing - the Game-of-Life state-chart image.
# ------------------------------------------------------------------
$shp = $d.InlineShapes.Item(1)
$r = $shp.Range
$start = $r.Start

# Pull the full WordOpenXML package so we can read the exact markup of the
# run that hosts the drawing (namespaces/rsids intact), then isolate just
# that <w:r>...</w:r> fragment.
$full = $r.WordOpenXML
$runStartMarker = "<w:r "
$drawMarker = "<w:drawing>"
$drawIdx = $full.IndexOf($drawMarker)
if ($drawIdx -lt 0) { throw "drawing marker not found in shape range XML" }
$runStart = $full.LastIndexOf($runStartMarker, $drawIdx)
if ($runStart -lt 0) { throw "enclosing run start not found" }
$runEndMarker = "</w:r>"
$runEnd = $full.IndexOf($runEndMarker, $drawIdx) + $runEndMarker.Length
$runXml = $full.Substring($runStart, $runEnd - $runStart)

# ------------------------------------------------------------------
# 1) Resize the drawing: wp:extent cx/cy shrink a bit.
# ------------------------------------------------------------------
$oldExtent = '<wp:extent cx="4933950" cy="3152775"/>'
$newExtent = '<wp:extent cx="4429125" cy="3009900"/>'
if ($runXml.IndexOf($oldExtent) -lt 0) { throw "wp:extent anchor not found" }
$runXml = $runXml.Replace($oldExtent, $newExtent)

# ------------------------------------------------------------------
# 2) Renumber the drawing's docPr id from 2 to 1.
# ------------------------------------------------------------------
$oldDocPr = '<wp:docPr id="2" name="Object 1"/>'
$newDocPr = '<wp:docPr id="1" name="Object 1"/>'
if ($runXml.IndexOf($oldDocPr) -lt 0) { throw "wp:docPr anchor not found" }
$runXml = $runXml.Replace($oldDocPr, $newDocPr)

# ------------------------------------------------------------------
# 3) Add a new dashed connector shape right after the first connector
#    ("Straight Connector 3") inside the locked canvas.
# ------------------------------------------------------------------
$firstConnector = '<a:cxnSp><a:nvCxnSpPr><a:cNvPr id="4" name="Straight Connector 3"/><a:cNvCxnSpPr/></a:nvCxnSpPr><a:spPr><a:xfrm><a:off x="3657600" y="3276600"/><a:ext cx="1066800" cy="0"/></a:xfrm><a:prstGeom prst="line"><a:avLst/></a:prstGeom><a:ln w="25400"/></a:spPr><a:style><a:lnRef idx="1"><a:schemeClr val="accent1"/></a:lnRef><a:fillRef idx="0"><a:schemeClr val="accent1"/></a:fillRef><a:effectRef idx="0"><a:schemeClr val="accent1"/></a:effectRef><a:fontRef idx="minor"><a:schemeClr val="tx1"/></a:fontRef></a:style></a:cxnSp>'

$newConnector = '<a:cxnSp><a:nvCxnSpPr><a:cNvPr id="3" name="Straight Connector 2"/><a:cNvCxnSpPr/></a:nvCxnSpPr><a:spPr><a:xfrm><a:off x="5867400" y="533400"/><a:ext cx="76200" cy="6172200"/></a:xfrm><a:prstGeom prst="line"><a:avLst/></a:prstGeom><a:ln><a:prstDash val="dash"/></a:ln></a:spPr><a:style><a:lnRef idx="2"><a:schemeClr val="accent1"/></a:lnRef><a:fillRef idx="0"><a:schemeClr val="accent1"/></a:fillRef><a:effectRef idx="1"><a:schemeClr val="accent1"/></a:effectRef><a:fontRef idx="minor"><a:schemeClr val="tx1"/></a:fontRef></a:style></a:cxnSp>'

if ($runXml.IndexOf($firstConnector) -lt 0) { throw "first connector anchor not found" }
$runXml = $runXml.Replace($firstConnector, $firstConnector + $newConnector)

# ------------------------------------------------------------------
# Remove the old drawing run and splice in the modified one at the same spot.
# ------------------------------------------------------------------
$beforeCount = $d.InlineShapes.Count
$shp.Delete()
if ($d.InlineShapes.Count -ne ($beforeCount - 1)) { throw "shape delete did not reduce InlineShapes.Count as expected" }

$pkg = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing"><w:body><w:p>' + $runXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target = $d.Range($start, $start)
$target.InsertXML($pkg)

if ($d.InlineShapes.Count -ne $beforeCount) { throw "re-inserted drawing did not restore InlineShapes.Count" }

Write-Host "Inline shapes after edit:" $d.InlineShapes.Count
Write-Host "Width:" $d.InlineShapes.Item(1).Width "Height:" $d.InlineShapes.Item(1).Height
